# April 2025 update without ECB
# Updates a handful of existing mean_fcast_sce values (tiny float-precision
# refresh), fills in B/C forecast values for Oct-2023..May-2024 rows, and
# appends six new monthly date rows (Sep-2024..Feb-2025) at the bottom of
# the series.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Small precision corrections on existing mean forecast values -----
$ws.Range("B48").Value  = 2.8627888042811942
$ws.Range("B52").Value  = 4.4007252819400593
$ws.Range("B68").Value  = 4.7951020103101882

# --- 2. Populate previously-blank B/C cells for rows 131-138 --------------
$ws.Range("B131").Value = 5.7662459745181556
$ws.Range("C131").Value = 5

$ws.Range("B132").Value = 4.6668984303790388
$ws.Range("C132").Value = 4

$ws.Range("B133").Value = 4.2313969537840279
$ws.Range("C133").Value = 4

$ws.Range("B134").Value = 5.6544956428051236
$ws.Range("C134").Value = 4

$ws.Range("B135").Value = 6.7576008482305978
$ws.Range("C135").Value = 5

$ws.Range("B136").Value = 3.6492519978627875
$ws.Range("C136").Value = 4

$ws.Range("B137").Value = 4.5776816546033379
$ws.Range("C137").Value = 4

$ws.Range("B138").Value = 5.0291237698224931
$ws.Range("C138").Value = 4

# --- 3. Append six new monthly date rows (142-147), B/C left blank --------
$ws.Range("A142").Value = 45536
$ws.Range("A143").Value = 45566
$ws.Range("A144").Value = 45597
$ws.Range("A145").Value = 45627
$ws.Range("A146").Value = 45658
$ws.Range("A147").Value = 45689

# carry the date-number formatting / B & C number formatting down from the
# preceding rows so the new rows match the existing style (s="1"/s="2")
$ws.Range("A141").Copy()
$ws.Range("A142:A147").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("B141:C141").Copy()
$ws.Range("B142:C147").PasteSpecial(-4122) # xlPasteFormats
